$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Change every yellow-highlighted paragraph to green highlight.
#    (All "yellow" highlight runs in the doc become "green" per diff.)
#    wdBrightGreen = 4 serializes to OOXML w:highlight w:val="green".
#    wdYellow      = 7 serializes to OOXML w:highlight w:val="yellow".
#    Using Paragraph.Range.Font.HighlightColorIndex also re-colors the
#    paragraph mark's own rPr (w:pPr/w:rPr/w:highlight), matching the
#    target XML exactly.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.HighlightColorIndex -eq 7) {
        $p.Range.Font.HighlightColorIndex = 4
    }
}

# ---------------------------------------------------------------------
# 2) "Seri" + _GoBack bookmark + "ja" -> single run "Serija".
#    A find/replace across the split re-merges the text into one run
#    and drops the bookmark that sat between the two runs.
#    MatchCase=$true so this does not also touch "...Karoserija..."
#    (which contains the case-insensitive substring "serija").
# ---------------------------------------------------------------------
$range = $d.Content
$range.Find.Execute("Serija", $true, $false, $false, $false, $false, $true, 1, $false, "Serija", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Re-create the _GoBack bookmark inside the „chat“ sustav bullet,
#    splitting "„chat“ sustav" into "„" + bookmark + "chat“ sustav".
# ---------------------------------------------------------------------
$quoteOpen = [char]0x201E
$quoteClose = [char]0x201C
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $txt = $p.Range.Text
    if ($txt -like "*${quoteOpen}chat${quoteClose} sustav*" -and $txt -notlike "*bez*" -and $txt -notlike "*nije*") {
        $r = $p.Range
        $r.Find.Execute($quoteOpen + "chat", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
        $splitPos = $r.Start + 1
        $bmRange = $d.Range($splitPos, $splitPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)
        break
    }
}
